$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.569.05'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '1.828.34'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.30'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3980'
$ws.Range('E8').Value = '  +5.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07816'
$ws.Range('E9').Value = '  +4.44%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.09'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.117'
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.340'
$ws.Range('E12').Value = '  +1.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.05'
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.002'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.566'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').Value = '1.828.54'
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.29'
$ws.Range('E17').Value = '  +4.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001093'
$ws.Range('E18').Value = '  +2.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06557'
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.81'
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('D23').Value = '28.588.24'
$ws.Range('E23').Value = '  +1.30%  '
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.241'
$ws.Range('E25').Value = '  +7.43%  '
$ws.Range('E26').Value = '  +1.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.12'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').Value = '2.040.42'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.414'
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.38'
$ws.Range('E30').Value = '  +2.73%  '
$ws.Range('E31').Value = '  +3.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1122'
$ws.Range('E32').Value = '  +2.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.744'
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.652'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07313'
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2265'
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.994'
$ws.Range('E37').Value = '  +6.37%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02349'
$ws.Range('E38').Value = '  +2.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.214'
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.39'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6295'
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.394'
$ws.Range('E44').Value = '  -3.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.49'
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5925'
$ws.Range('E46').Value = '  +2.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.712'
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.53'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.996'
$ws.Range('E49').Value = '  +3.55%  '
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('E51').Value = '  +1.77%  '
